$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: replace old "PNR Information" column with new, more
# granular PNR-related headers (Flight Number / PNR Number / PNR Validity /
# Origin / Destination).
$ws.Range("D1").Value = "Flight Number"
$ws.Range("E1").Value = "PNR Number"
$ws.Range("F1").Value = "PNR Validity"
$ws.Range("G1").Value = "Origin"
$ws.Range("H1").Value = "Destination"

# Fill in the corresponding data for the existing data row (row 2).
$ws.Range("D2").Value = "FR 202"
$ws.Range("E2").Value = "OTME5P"
$ws.Range("F2").Value = "Fri 06 Dec 2019"
$ws.Range("G2").Value = "DUBLIN"
$ws.Range("H2").Value = "LONDON STANSTED"

# Add a new test case row (row 3) for the booking-failed scenario.
$ws.Range("A3").Value = 2.0
$ws.Range("B3").Value = $ws.Range("B2").Value2
$ws.Range("C3").Value = $ws.Range("C2").Value2
$ws.Range("D3").Value = "AL 1235"
$ws.Range("E3").Value = "Booking Failed"

# PNR Validity / Origin / Destination are left blank for this row, but the
# source workbook stores them as explicit empty-text cells (shared-string
# entry), not truly-blank cells, so write them as empty text instead of
# leaving the cells untouched.
$ws.Range("F3").Formula = "'"
$ws.Range("F3").ClearFormats()
$ws.Range("G3").Formula = "'"
$ws.Range("G3").ClearFormats()
$ws.Range("H3").Formula = "'"
$ws.Range("H3").ClearFormats()
